$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Insert a new header column "Unnamed: 0.1.1" before the Date column; this
# pushes the existing Date/Grade headers (and their data) one column to the
# right (D->E, E->F).
$ws.Columns("D:D").Insert()

# New header cell needs the same bold/border/centered style as the other
# header cells.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Unnamed: 0.1.1"

# --- New "Unnamed: 0.1.1" data column (C) value for existing row 2 ---
$ws.Range("C2").Value = 0

# --- New "Unnamed: 0.1" data column (B) value for row 3 ---
$ws.Range("B3").Value = 1

# --- New row 4: third game record ---
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2

$ws.Range("E4").Value = "Sat Jan 18 19:20:10 2020"
$ws.Range("F4").Value = 80
